# Natmi following Dr Hou advice
# Recomputed the Col1a2 -> Cd93 ligand-receptor table: min expressing-cell
# threshold changed (most "1 cell" counts become "3 cells"), all derived
# expression/specificity statistics were updated accordingly, and three new
# sending-cluster rows (sCs -> ECs/FAPs/M2/sCs) were appended so every
# sending cluster x target cluster combination (ECs/FAPs/M2/sCs) is present.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Col1a2"
$ws.Range("C2").Value = "Cd93"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 6.423576
$ws.Range("H2").Value = 19.270728
$ws.Range("I2").Value = 0.001681024218962088
$ws.Range("J2").Value = 0.001681024218962088
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 135.955556
$ws.Range("N2").Value = 407.866668
$ws.Range("O2").Value = 0.6947679994035034
$ws.Range("P2").Value = 0.6947679994035034
$ws.Range("Q2").Value = 873.320846588256
$ws.Range("R2").Value = 7859.887619294303
$ws.Range("S2").Value = 0.001167921833557127
$ws.Range("T2").Value = 0.001167921833557127

# Row 3
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Col1a2"
$ws.Range("C3").Value = "Cd93"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 6.423576
$ws.Range("H3").Value = 19.270728
$ws.Range("I3").Value = 0.001681024218962088
$ws.Range("J3").Value = 0.001681024218962088
$ws.Range("K3").Value = 2
$ws.Range("L3").Value = 0.6666666666666666
$ws.Range("M3").Value = 0.449122
$ws.Range("N3").Value = 1.347366
$ws.Range("O3").Value = 0.002295129398228494
$ws.Range("P3").Value = 0.002295129398228494
$ws.Range("Q3").Value = 2.884969300272
$ws.Range("R3").Value = 25.964723702448
$ws.Range("S3").Value = 0.000003858168104073982
$ws.Range("T3").Value = 0.000003858168104073981

# Row 4
$ws.Range("A4").Value = "ECs"
$ws.Range("B4").Value = "Col1a2"
$ws.Range("C4").Value = "Cd93"
$ws.Range("D4").Value = "M2"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 6.423576
$ws.Range("H4").Value = 19.270728
$ws.Range("I4").Value = 0.001681024218962088
$ws.Range("J4").Value = 0.001681024218962088
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 56.38366533333333
$ws.Range("N4").Value = 169.150996
$ws.Range("O4").Value = 0.2881350899898248
$ws.Range("P4").Value = 0.2881350899898248
$ws.Range("Q4").Value = 362.184759427232
$ws.Range("R4").Value = 3259.662834845088
$ws.Range("S4").Value = 0.0004843620646057163
$ws.Range("T4").Value = 0.0004843620646057162

# Row 5
$ws.Range("A5").Value = "ECs"
$ws.Range("B5").Value = "Col1a2"
$ws.Range("C5").Value = "Cd93"
$ws.Range("D5").Value = "sCs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 6.423576
$ws.Range("H5").Value = 19.270728
$ws.Range("I5").Value = 0.001681024218962088
$ws.Range("J5").Value = 0.001681024218962088
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 2.896484
$ws.Range("N5").Value = 8.689452
$ws.Range("O5").Value = 0.01480178120844327
$ws.Range("P5").Value = 0.01480178120844327
$ws.Range("Q5").Value = 18.605785106784
$ws.Range("R5").Value = 167.452065961056
$ws.Range("S5").Value = 0.00002488215269517107
$ws.Range("T5").Value = 0.00002488215269517107

# Row 6
$ws.Range("A6").Value = "FAPs"
$ws.Range("B6").Value = "Col1a2"
$ws.Range("C6").Value = "Cd93"
$ws.Range("D6").Value = "ECs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 3580.644531333333
$ws.Range("H6").Value = 10741.933594
$ws.Range("I6").Value = 0.9370403925578976
$ws.Range("J6").Value = 0.9370403925578976
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 135.955556
$ws.Range("N6").Value = 407.866668
$ws.Range("O6").Value = 0.6947679994035034
$ws.Range("P6").Value = 0.6947679994035034
$ws.Range("Q6").Value = 486808.5180957828
$ws.Range("R6").Value = 4381276.662862045
$ws.Range("S6").Value = 0.651025678897724
$ws.Range("T6").Value = 0.651025678897724

# Row 7
$ws.Range("A7").Value = "FAPs"
$ws.Range("B7").Value = "Col1a2"
$ws.Range("C7").Value = "Cd93"
$ws.Range("D7").Value = "FAPs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 3580.644531333333
$ws.Range("H7").Value = 10741.933594
$ws.Range("I7").Value = 0.9370403925578976
$ws.Range("J7").Value = 0.9370403925578976
$ws.Range("K7").Value = 2
$ws.Range("L7").Value = 0.6666666666666666
$ws.Range("M7").Value = 0.449122
$ws.Range("N7").Value = 1.347366
$ws.Range("O7").Value = 0.002295129398228494
$ws.Range("P7").Value = 0.002295129398228494
$ws.Range("Q7").Value = 1608.146233201489
$ws.Range("R7").Value = 14473.3160988134
$ws.Range("S7").Value = 0.002150628952287199
$ws.Range("T7").Value = 0.002150628952287199

# Row 8
$ws.Range("A8").Value = "FAPs"
$ws.Range("B8").Value = "Col1a2"
$ws.Range("C8").Value = "Cd93"
$ws.Range("D8").Value = "M2"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 3580.644531333333
$ws.Range("H8").Value = 10741.933594
$ws.Range("I8").Value = 0.9370403925578976
$ws.Range("J8").Value = 0.9370403925578976
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 56.38366533333333
$ws.Range("N8").Value = 169.150996
$ws.Range("O8").Value = 0.2881350899898248
$ws.Range("P8").Value = 0.2881350899898248
$ws.Range("Q8").Value = 201889.8629323288
$ws.Range("R8").Value = 1817008.766390959
$ws.Range("S8").Value = 0.2699942178337706
$ws.Range("T8").Value = 0.2699942178337706

# Row 9
$ws.Range("A9").Value = "FAPs"
$ws.Range("B9").Value = "Col1a2"
$ws.Range("C9").Value = "Cd93"
$ws.Range("D9").Value = "sCs"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 3580.644531333333
$ws.Range("H9").Value = 10741.933594
$ws.Range("I9").Value = 0.9370403925578976
$ws.Range("J9").Value = 0.9370403925578976
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 2.896484
$ws.Range("N9").Value = 8.689452
$ws.Range("O9").Value = 0.01480178120844327
$ws.Range("P9").Value = 0.01480178120844327
$ws.Range("Q9").Value = 10371.2795946945
$ws.Range("R9").Value = 93341.51635225049
$ws.Range("S9").Value = 0.0138698668741158
$ws.Range("T9").Value = 0.0138698668741158

# Row 10
$ws.Range("A10").Value = "M2"
$ws.Range("B10").Value = "Col1a2"
$ws.Range("C10").Value = "Cd93"
$ws.Range("D10").Value = "ECs"
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 0.9157713333333334
$ws.Range("H10").Value = 2.747314
$ws.Range("I10").Value = 0.0002396537054071653
$ws.Range("J10").Value = 0.0002396537054071653
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 135.955556
$ws.Range("N10").Value = 407.866668
$ws.Range("O10").Value = 0.6947679994035034
$ws.Range("P10").Value = 0.6947679994035034
$ws.Range("Q10").Value = 124.5042007921947
$ws.Range("R10").Value = 1120.537807129752
$ws.Range("S10").Value = 0.0001665037254553728
$ws.Range("T10").Value = 0.0001665037254553728

# Row 11
$ws.Range("A11").Value = "M2"
$ws.Range("B11").Value = "Col1a2"
$ws.Range("C11").Value = "Cd93"
$ws.Range("D11").Value = "FAPs"
$ws.Range("E11").Value = 3
$ws.Range("F11").Value = 1
$ws.Range("G11").Value = 0.9157713333333334
$ws.Range("H11").Value = 2.747314
$ws.Range("I11").Value = 0.0002396537054071653
$ws.Range("J11").Value = 0.0002396537054071653
$ws.Range("K11").Value = 2
$ws.Range("L11").Value = 0.6666666666666666
$ws.Range("M11").Value = 0.449122
$ws.Range("N11").Value = 1.347366
$ws.Range("O11").Value = 0.002295129398228494
$ws.Range("P11").Value = 0.002295129398228494
$ws.Range("Q11").Value = 0.4112930527693334
$ws.Range("R11").Value = 3.701637474924
$ws.Range("S11").Value = 0.0000005500362646743759
$ws.Range("T11").Value = 0.0000005500362646743759

# Row 12
$ws.Range("A12").Value = "M2"
$ws.Range("B12").Value = "Col1a2"
$ws.Range("C12").Value = "Cd93"
$ws.Range("D12").Value = "M2"
$ws.Range("E12").Value = 3
$ws.Range("F12").Value = 1
$ws.Range("G12").Value = 0.9157713333333334
$ws.Range("H12").Value = 2.747314
$ws.Range("I12").Value = 0.0002396537054071653
$ws.Range("J12").Value = 0.0002396537054071653
$ws.Range("K12").Value = 3
$ws.Range("L12").Value = 1
$ws.Range("M12").Value = 56.38366533333333
$ws.Range("N12").Value = 169.150996
$ws.Range("O12").Value = 0.2881350899898248
$ws.Range("P12").Value = 0.2881350899898248
$ws.Range("Q12").Value = 51.63454438052711
$ws.Range("R12").Value = 464.710899424744
$ws.Range("S12").Value = 0.00006905264197388853
$ws.Range("T12").Value = 0.00006905264197388853

# Row 13
$ws.Range("A13").Value = "M2"
$ws.Range("B13").Value = "Col1a2"
$ws.Range("C13").Value = "Cd93"
$ws.Range("D13").Value = "sCs"
$ws.Range("E13").Value = 3
$ws.Range("F13").Value = 1
$ws.Range("G13").Value = 0.9157713333333334
$ws.Range("H13").Value = 2.747314
$ws.Range("I13").Value = 0.0002396537054071653
$ws.Range("J13").Value = 0.0002396537054071653
$ws.Range("K13").Value = 3
$ws.Range("L13").Value = 1
$ws.Range("M13").Value = 2.896484
$ws.Range("N13").Value = 8.689452
$ws.Range("O13").Value = 0.01480178120844327
$ws.Range("P13").Value = 0.01480178120844327
$ws.Range("Q13").Value = 2.652517014658666
$ws.Range("R13").Value = 23.872653131928
$ws.Range("S13").Value = 0.000003547301713229579
$ws.Range("T13").Value = 0.000003547301713229579

# Row 14
$ws.Range("A14").Value = "sCs"
$ws.Range("B14").Value = "Col1a2"
$ws.Range("C14").Value = "Cd93"
$ws.Range("D14").Value = "ECs"
$ws.Range("E14").Value = 3
$ws.Range("F14").Value = 1
$ws.Range("G14").Value = 233.243637
$ws.Range("H14").Value = 699.7309110000001
$ws.Range("I14").Value = 0.0610389295177331
$ws.Range("J14").Value = 0.06103892951773311
$ws.Range("K14").Value = 3
$ws.Range("L14").Value = 1
$ws.Range("M14").Value = 135.955556
$ws.Range("N14").Value = 407.866668
$ws.Range("O14").Value = 0.6947679994035034
$ws.Range("P14").Value = 0.6947679994035034
$ws.Range("Q14").Value = 31710.76835179718
$ws.Range("R14").Value = 285396.9151661746
$ws.Range("S14").Value = 0.04240789494676687
$ws.Range("T14").Value = 0.04240789494676688

# Row 15
$ws.Range("A15").Value = "sCs"
$ws.Range("B15").Value = "Col1a2"
$ws.Range("C15").Value = "Cd93"
$ws.Range("D15").Value = "FAPs"
$ws.Range("E15").Value = 3
$ws.Range("F15").Value = 1
$ws.Range("G15").Value = 233.243637
$ws.Range("H15").Value = 699.7309110000001
$ws.Range("I15").Value = 0.0610389295177331
$ws.Range("J15").Value = 0.06103892951773311
$ws.Range("K15").Value = 2
$ws.Range("L15").Value = 0.6666666666666666
$ws.Range("M15").Value = 0.449122
$ws.Range("N15").Value = 1.347366
$ws.Range("O15").Value = 0.002295129398228494
$ws.Range("P15").Value = 0.002295129398228494
$ws.Range("Q15").Value = 104.754848736714
$ws.Range("R15").Value = 942.7936386304261
$ws.Range("S15").Value = 0.0001400922415725462
$ws.Range("T15").Value = 0.0001400922415725462

# Row 16
$ws.Range("A16").Value = "sCs"
$ws.Range("B16").Value = "Col1a2"
$ws.Range("C16").Value = "Cd93"
$ws.Range("D16").Value = "M2"
$ws.Range("E16").Value = 3
$ws.Range("F16").Value = 1
$ws.Range("G16").Value = 233.243637
$ws.Range("H16").Value = 699.7309110000001
$ws.Range("I16").Value = 0.0610389295177331
$ws.Range("J16").Value = 0.06103892951773311
$ws.Range("K16").Value = 3
$ws.Range("L16").Value = 1
$ws.Range("M16").Value = 56.38366533333333
$ws.Range("N16").Value = 169.150996
$ws.Range("O16").Value = 0.2881350899898248
$ws.Range("P16").Value = 0.2881350899898248
$ws.Range("Q16").Value = 13151.13116973749
$ws.Range("R16").Value = 118360.1805276374
$ws.Range("S16").Value = 0.0175874574494746
$ws.Range("T16").Value = 0.0175874574494746

# Row 17
$ws.Range("A17").Value = "sCs"
$ws.Range("B17").Value = "Col1a2"
$ws.Range("C17").Value = "Cd93"
$ws.Range("D17").Value = "sCs"
$ws.Range("E17").Value = 3
$ws.Range("F17").Value = 1
$ws.Range("G17").Value = 233.243637
$ws.Range("H17").Value = 699.7309110000001
$ws.Range("I17").Value = 0.0610389295177331
$ws.Range("J17").Value = 0.06103892951773311
$ws.Range("K17").Value = 3
$ws.Range("L17").Value = 1
$ws.Range("M17").Value = 2.896484
$ws.Range("N17").Value = 8.689452
$ws.Range("O17").Value = 0.01480178120844327
$ws.Range("P17").Value = 0.01480178120844327
$ws.Range("Q17").Value = 675.586462672308
$ws.Range("R17").Value = 6080.278164050773
$ws.Range("S17").Value = 0.0009034848799190752
$ws.Range("T17").Value = 0.0009034848799190753

